$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: update to new publish date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; now populated
$ws.Range("B9").Value = "Alvearie Team"

# Remove the duplicated "Contact" row (rows 10 & 11 were identical)
$ws.Rows.Item(11).Delete()

# The remaining former "Contact" row (row 10) becomes "Jurisdiction"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

Write-Host "Done"
